$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 25, shifting existing rows 25-55 down to 26-56
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new weekly price observation
$ws.Cells.Item(25, 1).Value = 4
$ws.Cells.Item(25, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(25, 3).Value = "Los Lagos"
$ws.Cells.Item(25, 4).Value = 44967
$ws.Cells.Item(25, 5).Value = 10
$ws.Cells.Item(25, 6).Value = "Fruta"
$ws.Cells.Item(25, 7).Value = 100101
$ws.Cells.Item(25, 8).Value = "Berries"
$ws.Cells.Item(25, 9).Value = 100101001
$ws.Cells.Item(25, 10).Value = "Arándano (blue)"
$ws.Cells.Item(25, 11).Value = "Sin especificar"
$ws.Cells.Item(25, 12).Value = "Primera"
$ws.Cells.Item(25, 13).Value = 300
$ws.Cells.Item(25, 14).Value = 2000
$ws.Cells.Item(25, 15).Value = 2200
$ws.Cells.Item(25, 16).Value = 2100
$ws.Cells.Item(25, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(25, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(25, 19).Value = 1050
$ws.Cells.Item(25, 20).Value = 2
